# Update lattice multiplication exercise table cells to new values.
# Each cell holds 5 lines of text joined by manual line breaks; inside
# Range.Text a <w:br/> shows up as a vertical-tab character (char code 11).
$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

function Set-LatticeCell {
    param($Table, $Row, $Col, $ExpectedFirstLine, $Lines)
    $vt = [char]11
    $cell = $Table.Cell($Row, $Col)
    $current = $cell.Range.Text
    if ($current.IndexOf($ExpectedFirstLine) -ne 0) {
        Write-Host "WARNING: cell ($Row,$Col) did not start with $ExpectedFirstLine (got $current)"
    }
    $cell.Range.Text = [string]::Join($vt, $Lines)
}

# Row 1, Col 1: "48 x 94" -> "16 x 45"
Set-LatticeCell $table 1 1 "48 x 94" @("16 x 45", "  4    5", "  ----", "1|    |", "6|    |")
# Row 1, Col 2: "25 x 65" -> "10 x 15"
Set-LatticeCell $table 1 2 "25 x 65" @("10 x 15", "  1    5", "  ----", "1|    |", "0|    |")
# Row 1, Col 3: "57 x 25" -> "19 x 34"
Set-LatticeCell $table 1 3 "57 x 25" @("19 x 34", "  3    4", "  ----", "1|    |", "9|    |")
# Row 2, Col 1: "21 x 84" -> "71 x 93"
Set-LatticeCell $table 2 1 "21 x 84" @("71 x 93", "  9    3", "  ----", "7|    |", "1|    |")
# Row 2, Col 2: "83 x 59" -> "91 x 21"
Set-LatticeCell $table 2 2 "83 x 59" @("91 x 21", "  2    1", "  ----", "9|    |", "1|    |")
# Row 2, Col 3: "93 x 56" -> "24 x 40"
Set-LatticeCell $table 2 3 "93 x 56" @("24 x 40", "  4    0", "  ----", "2|    |", "4|    |")
# Row 3, Col 1: "86 x 56" -> "47 x 47"
Set-LatticeCell $table 3 1 "86 x 56" @("47 x 47", "  4    7", "  ----", "4|    |", "7|    |")
# Row 3, Col 2: "63 x 54" -> "60 x 91"
Set-LatticeCell $table 3 2 "63 x 54" @("60 x 91", "  9    1", "  ----", "6|    |", "0|    |")
# Row 3, Col 3: "64 x 39" -> "59 x 68"
Set-LatticeCell $table 3 3 "64 x 39" @("59 x 68", "  6    8", "  ----", "5|    |", "9|    |")
# Row 4, Col 1: "47 x 54" -> "87 x 98"
Set-LatticeCell $table 4 1 "47 x 54" @("87 x 98", "  9    8", "  ----", "8|    |", "7|    |")
# Row 4, Col 2: "95 x 19" -> "59 x 65"
Set-LatticeCell $table 4 2 "95 x 19" @("59 x 65", "  6    5", "  ----", "5|    |", "9|    |")
# Row 4, Col 3: "19 x 56" -> "59 x 78"
Set-LatticeCell $table 4 3 "19 x 56" @("59 x 78", "  7    8", "  ----", "5|    |", "9|    |")
# Row 5, Col 1: "86 x 35" -> "98 x 41"
Set-LatticeCell $table 5 1 "86 x 35" @("98 x 41", "  4    1", "  ----", "9|    |", "8|    |")
# Row 5, Col 2: "91 x 11" -> "51 x 70"
Set-LatticeCell $table 5 2 "91 x 11" @("51 x 70", "  7    0", "  ----", "5|    |", "1|    |")
# Row 5, Col 3: "35 x 59" -> "56 x 46"
Set-LatticeCell $table 5 3 "35 x 59" @("56 x 46", "  4    6", "  ----", "5|    |", "6|    |")

Write-Host "Done updating lattice multiplication table."
